$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 3893.8333
$ws.Range("I32").Value = 3655.3333
$ws.Range("J32").Value = 4132.3335
$ws.Range("K32").Value = 3655.3333
$ws.Range("L32").Value = 4132.3335
$ws.Range("M32").Value = -3329.3333
$ws.Range("N32").Value = -4784.3335
$ws.Range("H52").Value = 300
$ws.Range("I52").Value = 300
$ws.Range("K52").Value = 900
$ws.Range("M52").Value = -740
$ws.Range("H58").Value = 959.1539
$ws.Range("J58").Value = 2624.5
$ws.Range("L58").Value = 7873.5
$ws.Range("N58").Value = -8173.5
$ws.Range("H70").Value = 2050.8125
$ws.Range("I70").Value = 1413.8572
$ws.Range("J70").Value = 2546.2222
$ws.Range("K70").Value = 4241.571599999999
$ws.Range("L70").Value = 7638.6666
$ws.Range("M70").Value = -3971.571599999999
$ws.Range("N70").Value = -8178.6666
$ws.Range("H73").Value = 2050.8125
$ws.Range("I73").Value = 1413.8572
$ws.Range("J73").Value = 2546.2222
$ws.Range("K73").Value = 4241.571599999999
$ws.Range("L73").Value = 7638.6666
$ws.Range("M73").Value = -3305.571599999999
$ws.Range("N73").Value = -9510.6666
$ws.Range("H98").Value = 3747.4
$ws.Range("I98").Value = 2947.3333
$ws.Range("J98").Value = 4947.5
$ws.Range("K98").Value = 2947.3333
$ws.Range("L98").Value = 4947.5
$ws.Range("M98").Value = -1449.3333
$ws.Range("N98").Value = -7943.5
$ws.Range("H122").Value = 3747.4
$ws.Range("I122").Value = 2947.3333
$ws.Range("J122").Value = 4947.5
$ws.Range("K122").Value = 8841.999899999999
$ws.Range("L122").Value = 14842.5
$ws.Range("M122").Value = -6391.999899999999
$ws.Range("N122").Value = -19742.5
$ws.Range("H132").Value = 3698.322
$ws.Range("I132").Value = 2149.6365
$ws.Range("K132").Value = 6448.9095
$ws.Range("M132").Value = -3918.9095
$ws.Range("H137").Value = 4340.3687
$ws.Range("I137").Value = 969
$ws.Range("K137").Value = 2907
$ws.Range("M137").Value = -357

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H53").Value = 0
$ws.Range("I53").Value = 0
$ws.Range("K53").Value = 0
$ws.Range("M53").ClearContents()
$ws.Range("H74").Value = 167481.83
$ws.Range("I74").Value = 167481.83
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 167481.83
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()
$ws.Range("M74").Value = -166607.83
$ws.Range("H77").Value = 167481.83
$ws.Range("I77").Value = 167481.83
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 837409.1499999999
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()
$ws.Range("M77").Value = -833041.1499999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 4557.75
$ws.Range("I105").Value = 4432.04
$ws.Range("K105").Value = 4432.04
$ws.Range("M105").Value = -2685.04
$ws.Range("H107").Value = 1438.75
$ws.Range("I107").Value = 1438.75
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 1438.75
$ws.Range("L107").Value = 0
$ws.Range("N107").ClearContents()
$ws.Range("M107").Value = 481.25
$ws.Range("H134").Value = 2189.375
$ws.Range("I134").Value = 1902.5264
$ws.Range("J134").Value = 3279.4
$ws.Range("K134").Value = 5707.5792
$ws.Range("L134").Value = 9838.200000000001
$ws.Range("M134").Value = -3172.5792
$ws.Range("N134").Value = -14908.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 87
$ws.Range("I7").Value = 63.57143
$ws.Range("J7").Value = 119.8
$ws.Range("K7").Value = 63.57143
$ws.Range("L7").Value = 119.8
$ws.Range("M7").Value = 49.42857
$ws.Range("N7").Value = -345.8
$ws.Range("H47").Value = 34750
$ws.Range("I47").Value = 33000
$ws.Range("K47").Value = 33000
$ws.Range("M47").Value = -32434
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("N63").ClearContents()
$ws.Range("L63").Value = 0
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("N66").ClearContents()
$ws.Range("L66").Value = 0
$ws.Range("H86").Value = 45257
$ws.Range("I86").Value = 62182.43
$ws.Range("J86").Value = 5764.3335
$ws.Range("K86").Value = 62182.43
$ws.Range("L86").Value = 5764.3335
$ws.Range("M86").Value = -61059.43
$ws.Range("N86").Value = -8010.3335
$ws.Range("H89").Value = 45257
$ws.Range("I89").Value = 62182.43
$ws.Range("J89").Value = 5764.3335
$ws.Range("K89").Value = 310912.15
$ws.Range("L89").Value = 28821.6675
$ws.Range("M89").Value = -305296.15
$ws.Range("N89").Value = -40053.6675
$ws.Range("H107").Value = 365.7143
$ws.Range("J107").Value = 533.3333
$ws.Range("L107").Value = 533.3333
$ws.Range("N107").Value = -4373.3333
$ws.Range("H122").Value = 3699.25
$ws.Range("I122").Value = 3398.5
$ws.Range("K122").Value = 10195.5
$ws.Range("M122").Value = -7745.5
$ws.Range("H132").Value = 6241
$ws.Range("I132").Value = 5301.375
$ws.Range("J132").Value = 9999.5
$ws.Range("K132").Value = 15904.125
$ws.Range("L132").Value = 29998.5
$ws.Range("M132").Value = -13374.125
$ws.Range("N132").Value = -35058.5
$ws.Range("H134").Value = 43122.68
$ws.Range("I134").Value = 48593.953
$ws.Range("J134").Value = 3000
$ws.Range("K134").Value = 145781.859
$ws.Range("L134").Value = 9000
$ws.Range("M134").Value = -143246.859
$ws.Range("N134").Value = -14070

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 2224.5
$ws.Range("I3").Value = 2224.5
$ws.Range("K3").Value = 6673.5
$ws.Range("M3").Value = -6561.5
$ws.Range("H13").Value = 61
$ws.Range("I13").Value = 13.2
$ws.Range("K13").Value = 39.59999999999999
$ws.Range("M13").Value = 128.4
$ws.Range("H46").Value = 777
$ws.Range("I46").Value = 777
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 2331
$ws.Range("L46").Value = 0
$ws.Range("N46").ClearContents()
$ws.Range("M46").Value = -2240
$ws.Range("H51").Value = 833.3333
$ws.Range("I51").Value = 833.3333
$ws.Range("K51").Value = 2499.9999
$ws.Range("M51").Value = -2039.9999
$ws.Range("H54").Value = 0
$ws.Range("I54").Value = 0
$ws.Range("K54").Value = 0
$ws.Range("M54").ClearContents()
$ws.Range("H113").Value = 1337.8379
$ws.Range("I113").Value = 1054.7778
$ws.Range("J113").Value = 1428.8214
$ws.Range("K113").Value = 3164.3334
$ws.Range("L113").Value = 4286.4642
$ws.Range("M113").Value = -994.3334000000004
$ws.Range("N113").Value = -8626.4642
$ws.Range("H118").Value = 5536.25
$ws.Range("I118").Value = 5536.25
$ws.Range("K118").Value = 16608.75
$ws.Range("M118").Value = -15365.75
$ws.Range("H131").Value = 4008507
$ws.Range("J131").Value = 4771373.5
$ws.Range("L131").Value = 14314120.5
$ws.Range("N131").Value = -14324200.5
$ws.Range("H132").Value = 1049.5454
$ws.Range("I132").Value = 820.8570999999999
$ws.Range("K132").Value = 7387.7139
$ws.Range("M132").Value = -4857.7139
$ws.Range("H134").Value = 2660.318
$ws.Range("I134").Value = 1825.1177
$ws.Range("K134").Value = 5475.3531
$ws.Range("M134").Value = -405.3531000000003
$ws.Range("H136").Value = 1749.5
$ws.Range("I136").Value = 1749.5
$ws.Range("K136").Value = 5248.5
$ws.Range("M136").Value = -148.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2268.6453
$ws.Range("I102").Value = 2082.88
$ws.Range("K102").Value = 2082.88
$ws.Range("M102").Value = -460.8800000000001
$ws.Range("H107").Value = 84707.336
$ws.Range("I107").Value = 200447.8
$ws.Range("J107").Value = 2035.5714
$ws.Range("K107").Value = 200447.8
$ws.Range("L107").Value = 2035.5714
$ws.Range("M107").Value = -198527.8
$ws.Range("N107").Value = -5875.5714
$ws.Range("H136").Value = 39851.5
$ws.Range("J136").Value = 39851.5
$ws.Range("L136").Value = 119554.5
$ws.Range("N136").Value = -124654.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 47629.074
$ws.Range("I132").Value = 66030.78999999999
$ws.Range("K132").Value = 198092.37
$ws.Range("M132").Value = -195562.37

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 9999
$ws.Range("J5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("N5").ClearContents()
$ws.Range("H122").Value = 3792.1667
$ws.Range("I122").Value = 4213.25
$ws.Range("K122").Value = 12639.75
$ws.Range("M122").Value = -10189.75
$ws.Range("H132").Value = 418143.28
$ws.Range("I132").Value = 486167.16
$ws.Range("K132").Value = 1458501.48
$ws.Range("M132").Value = -1455971.48
